# Auto-generated cell updates derived from the authoritative xml diff (cryptos.xlsx)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    # Force the value to be stored as literal text (matches the source
    # workbook's inlineStr cells) instead of letting Excel auto-coerce
    # numeric-looking strings (e.g. '1.00', '0.0000247') into numbers.
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range($addr).Style = "Normal"
}

Set-TextValue 'D2' '93.393.66'
Set-TextValue 'E2' '  -0.06%  '
Set-TextValue 'D3' '3.412.09'
Set-TextValue 'E3' '  +1.64%  '
Set-TextValue 'E4' '  +0.00%  '
Set-TextValue 'D5' '233.52'
Set-TextValue 'E5' '  -0.21%  '
Set-TextValue 'D6' '620.49'
Set-TextValue 'E6' '  -1.39%  '
Set-TextValue 'D7' '1.46'
Set-TextValue 'E7' '  +6.76%  '
Set-TextValue 'D8' '0.391'
Set-TextValue 'E8' '  +0.02%  '
Set-TextValue 'E9' '  -0.01%  '
Set-TextValue 'D10' '0.990'
Set-TextValue 'E10' '  +5.39%  '
Set-TextValue 'D11' '3.410.69'
Set-TextValue 'E11' '  +1.61%  '
Set-TextValue 'D12' '43.36'
Set-TextValue 'E12' '  +7.32%  '
Set-TextValue 'E13' '  +2.17%  '
Set-TextValue 'D14' '6.28'
Set-TextValue 'E14' '  +4.70%  '
Set-TextValue 'D15' '93.199.61'
Set-TextValue 'E15' '  -0.04%  '
Set-TextValue 'D16' '4.061.94'
Set-TextValue 'E16' '  +2.00%  '
Set-TextValue 'D17' '0.0000247'
Set-TextValue 'E17' '  +1.16%  '
Set-TextValue 'D18' '8.29'
Set-TextValue 'E18' '  +3.72%  '
Set-TextValue 'D19' '3.409.65'
Set-TextValue 'E19' '  +1.64%  '
Set-TextValue 'D20' '18.04'
Set-TextValue 'E20' '  +6.92%  '
Set-TextValue 'E21' '  +6.72%  '
Set-TextValue 'D22' '0.511'
Set-TextValue 'E22' '  +13.23%  '
Set-TextValue 'D23' '3.39'
Set-TextValue 'E23' '  +8.14%  '
Set-TextValue 'D24' '498.44'
Set-TextValue 'E24' '  +1.01%  '
Set-TextValue 'D25' '6.79'
Set-TextValue 'E25' '  +7.86%  '
Set-TextValue 'D26' '0.0000183'
Set-TextValue 'E26' '  -1.69%  '
Set-TextValue 'D27' '90.04'
Set-TextValue 'E27' '  +0.27%  '
Set-TextValue 'D28' '12.03'
Set-TextValue 'E28' '  +4.82%  '
Set-TextValue 'D29' '11.34'
Set-TextValue 'E29' '  +0.19%  '
Set-TextValue 'B30' 'Hedera'
Set-TextValue 'C30' 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue 'D30' '0.141'
Set-TextValue 'E30' '  +6.60%  '
Set-TextValue 'B31' 'Dai'
Set-TextValue 'C31' 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue 'D31' '1.00'
Set-TextValue 'E31' '  +0.01%  '
Set-TextValue 'D32' '2.72'
Set-TextValue 'E32' '  +2.63%  '
Set-TextValue 'D33' '1.01'
Set-TextValue 'E33' '  +0.70%  '
Set-TextValue 'E34' '  +2.23%  '
Set-TextValue 'D35' '0.551'
Set-TextValue 'E35' '  +4.56%  '
Set-TextValue 'D36' '28.92'
Set-TextValue 'E36' '  +0.88%  '
Set-TextValue 'D37' '558.25'
Set-TextValue 'E37' '  +6.85%  '
Set-TextValue 'D38' '7.49'
Set-TextValue 'E38' '  +0.08%  '
Set-TextValue 'B39' 'Fetch.AI'
Set-TextValue 'C39' 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue 'D39' '1.41'
Set-TextValue 'E39' '  +0.93%  '
Set-TextValue 'B40' 'USDe'
Set-TextValue 'C40' 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
Set-TextValue 'D40' '1.00'
Set-TextValue 'E40' '  -0.04%  '
Set-TextValue 'E41' '  +1.39%  '
Set-TextValue 'E42' '  +2.06%  '
Set-TextValue 'D43' '23.71'
Set-TextValue 'E43' '  -1.34%  '
Set-TextValue 'D44' '1.70'
Set-TextValue 'E44' '  +1.86%  '
Set-TextValue 'B45' 'VeChain'
Set-TextValue 'C45' 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue 'D45' '0.0418'
Set-TextValue 'E45' '  +6.29%  '
Set-TextValue 'B46' 'MantraDAO'
Set-TextValue 'C46' 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
Set-TextValue 'D46' '3.63'
Set-TextValue 'E46' '  +0.90%  '
Set-TextValue 'D47' '5.46'
Set-TextValue 'D48' '53.05'
Set-TextValue 'E48' '  +1.38%  '
Set-TextValue 'D49' '8.10'
Set-TextValue 'E49' '  +1.60%  '
Set-TextValue 'E50' '  -2.20%  '
Set-TextValue 'D51' '3.07'
Set-TextValue 'E51' '  -1.98%  '
